$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "27.439.38"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "1.743.97"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "321.64"
$ws.Range("E5").Value = "  -4.44%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.4195"
$ws.Range("E7").Value = "  -9.17%  "
$ws.Range("D8").Value = "0.3578"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("D9").Value = "45.40"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "0.07414"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("D11").Value = "1.110"
$ws.Range("E11").Value = "  -3.67%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "21.36"
$ws.Range("E13").Value = "  -4.58%  "
$ws.Range("D14").Value = "6.100"
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").Value = "7.178"
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("D16").Value = "1.744.94"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").Value = "0.00001062"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").Value = "88.01"
$ws.Range("E18").Value = "  +7.34%  "
$ws.Range("D19").Value = "0.06114"
$ws.Range("E19").Value = "  -9.08%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "16.82"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "6.087"
$ws.Range("E22").Value = "  -5.20%  "
$ws.Range("D23").Value = "0.5271"
$ws.Range("E23").Value = "  -5.26%  "
$ws.Range("D24").Value = "27.477.51"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").Value = "11.45"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("D26").Value = "2.336"
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("D27").Value = "20.35"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "152.85"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "2.367"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "1.942.04"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("D31").Value = "125.63"
$ws.Range("E31").Value = "  -5.87%  "
$ws.Range("D32").Value = "1.190"
$ws.Range("E32").Value = "  -5.43%  "
$ws.Range("D33").Value = "5.643"
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("D34").Value = "0.09113"
$ws.Range("E34").Value = "  -4.74%  "
$ws.Range("D35").Value = "3.630"
$ws.Range("E35").Value = "  -10.00%  "
$ws.Range("D36").Value = "12.58"
$ws.Range("E36").Value = "  +3.62%  "
$ws.Range("D37").Value = "0.02289"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").Value = "0.2132"
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("D39").Value = "5.068"
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("D40").Value = "0.06040"
$ws.Range("E40").Value = "  -5.16%  "
$ws.Range("D41").Value = "0.6369"
$ws.Range("E41").Value = "  -4.31%  "
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("D43").Value = "1.438"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "7.866"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("D46").Value = "13.74"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("D47").Value = "3.713"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").Value = "0.5838"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("D49").Value = "124.81"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("D50").Value = "1.940"
$ws.Range("E50").Value = "  -5.54%  "

Write-Host "Updated cryptos list"
